$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price cells (col D) get purely numeric-looking text (e.g. "213.43").
# Force those specific cells to Text format first so Excel keeps them as text
# (matching the original inlineStr typing) instead of auto-converting to numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "27.943.24"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.646.22"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "213.43"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D8").Value = "23.46"
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "1.880.25"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.651.00"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "65.64"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "27.964.66"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "231.66"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "7.66"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D22").Value = "10.66"
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("D23").Value = "4.39"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").Value = "152.23"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "15.76"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").Value = "1.440.28"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "0.889"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.935"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").Value = "0.558"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "69.17"
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "1.83"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("D46").Value = "5.43"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "1.788.93"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("E51").Value = "  -0.04%  "
